# First commit for 9th march 2017
# Adds 4 new rows (104-107) of daily-log data to the bottom of sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Populate new string values in the exact order needed so that the
# shared-strings table is appended in the same sequence as the source
# workbook (9th Mar,2017 / Webflow / Spring security / Cart
# implementation / pending / 360 minutes).
# ---------------------------------------------------------------------
$ws.Range("B104").Value = "9th Mar,2017"
$ws.Range("C104").Value = "Webflow "
$ws.Range("C105").Value = "Spring security"
$ws.Range("C107").Value = "Cart implementation"
$ws.Range("F107").Value = "pending"
$ws.Range("F104").Value = "360 minutes"

# ---------------------------------------------------------------------
# Row 104
# ---------------------------------------------------------------------
$ws.Range("A104").Value = 103
$ws.Range("D104").Value = "NA"
$ws.Range("E104").Value = "NA"
$ws.Range("G104").Value = "N"
$ws.Range("H104").Value = "NA"

# ---------------------------------------------------------------------
# Row 105
# ---------------------------------------------------------------------
$ws.Range("A105").Value = 104
$ws.Range("B105").Value = "9th Mar,2017"
$ws.Range("D105").Value = "NA"
$ws.Range("E105").Value = "NA"
$ws.Range("F105").Value = "360 minutes"
$ws.Range("G105").Value = "N"
$ws.Range("H105").Value = "NA"

# ---------------------------------------------------------------------
# Row 106
# ---------------------------------------------------------------------
$ws.Range("A106").Value = 105
$ws.Range("B106").Value = "9th Mar,2017"
$ws.Range("C106").Value = "Documentation"
$ws.Range("D106").Value = "NA"
$ws.Range("E106").Value = "NA"
$ws.Range("F106").Value = "30 minutes"

# D106 carries the "hyperlink-like" cell style (s="20") used elsewhere in
# the sheet for this column even though no live hyperlink is attached;
# copy that formatting from an existing cell that already uses it.
$ws.Range("D22").Copy() | Out-Null
$ws.Range("D106").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Row 107
# ---------------------------------------------------------------------
$ws.Range("A107").Value = 106
$ws.Range("B107").Value = "9th Mar,2017"

# ---------------------------------------------------------------------
# Row heights - the source content wraps onto two lines in these rows
# (matching every other multi-line row in the sheet).
# ---------------------------------------------------------------------
$ws.Rows.Item(104).RowHeight = 28.8
$ws.Rows.Item(105).RowHeight = 28.8
$ws.Rows.Item(106).RowHeight = 28.8
$ws.Rows.Item(107).RowHeight = 28.8

# ---------------------------------------------------------------------
# View state: scroll down and leave the selection on the last cell
# touched, matching where the author ended up after data entry.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 97
$ws.Range("F107").Select()
